$wb = $excel.ActiveWorkbook

# Sheet 1: 展览 (Exhibition)
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F2").Value = 10556
$ws1.Range("F6").Value = 292
$ws1.Range("F8").Value = 481
$ws1.Range("F9").Value = 790
$ws1.Range("F11").Value = 1217
$ws1.Range("F12").Value = 1116
$ws1.Range("F13").Value = 3269
$ws1.Range("F14").Value = 2428
$ws1.Range("F16").Value = 2194
$ws1.Range("F20").Value = 1601
$ws1.Range("F21").Value = 595
$ws1.Range("F23").Value = 261
$ws1.Range("F28").Value = 391
$ws1.Range("F30").Value = 79
$ws1.Range("F32").Value = 613
$ws1.Range("F33").Value = 37
$ws1.Range("F34").Value = 59
$ws1.Range("F35").Value = 286
$ws1.Range("F36").Value = 15
$ws1.Range("F38").Value = 540
$ws1.Range("F39").Value = 503
$ws1.Range("F40").Value = 1750
$ws1.Range("F41").Value = 154
$ws1.Range("F42").Value = 455
$ws1.Range("F43").Value = 61
$ws1.Range("F44").Value = 484
$ws1.Range("F45").Value = 1066

# Sheet 2: 演出 (Performance)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F4").Value = 52

# Sheet 4: 全部类型 (All types)
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F2").Value = 10556
$ws4.Range("F8").Value = 292
$ws4.Range("F10").Value = 481
$ws4.Range("F11").Value = 790
$ws4.Range("F12").Value = 1116
$ws4.Range("F13").Value = 3269
$ws4.Range("F14").Value = 2428
$ws4.Range("F15").Value = 2194
$ws4.Range("F16").Value = 1601
$ws4.Range("F17").Value = 595
$ws4.Range("F19").Value = 261
$ws4.Range("F24").Value = 391
$ws4.Range("F26").Value = 79
$ws4.Range("F28").Value = 613
$ws4.Range("F29").Value = 37
$ws4.Range("F30").Value = 52
$ws4.Range("F33").Value = 59
$ws4.Range("F34").Value = 286
$ws4.Range("F36").Value = 540
$ws4.Range("F38").Value = 503
$ws4.Range("F39").Value = 1750
$ws4.Range("F40").Value = 154
$ws4.Range("F44").Value = 455
$ws4.Range("F45").Value = 61
$ws4.Range("F46").Value = 484
$ws4.Range("F47").Value = 1066
